# Centralized Clearing workbook update:
# Insert a new "CL.NO.MV" column (column C) into both the CCSameDayInputter
# and CCNormalDayInputter sheets, pushing the existing C:F columns to D:G,
# then restore the view state (active sheet / selection) the author ended
# up with.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # CCSameDayInputter
$ws2 = $wb.Worksheets.Item(2)   # CCNormalDayInputter

# --- CCSameDayInputter: insert column C, label it, match column B's width ---
$ws1.Columns.Item(3).Insert()
$ws1.Cells.Item(1, 3).Value = "CL.NO.MV"
$ws1.Columns.Item(3).ColumnWidth = $ws1.Columns.Item(2).ColumnWidth

# --- CCNormalDayInputter: same column insert ---
$ws2.Columns.Item(3).Insert()
$ws2.Cells.Item(1, 3).Value = "CL.NO.MV"
$ws2.Columns.Item(3).ColumnWidth = $ws2.Columns.Item(2).ColumnWidth

# --- Restore selection/active-sheet state ---
# CCNormalDayInputter ends up with plain C1 selected (no longer the active tab)
$ws2.Range("C1").Select()
# CCSameDayInputter becomes the active tab, with G16 selected
$ws1.Range("G16").Select()
